# Update "Pais" (countries) data sheet and "Datos actualizados" timestamp,
# plus reorder Santo Tome y Principe / Surinam / Papua Nueva Guinea rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 10:22"

# --- Suiza (row 20): only new deaths (F) changes ---
$ws.Range("F20").Value = 167

# --- Austria (row 30) ---
$ws.Range("B30").Value = 15452
$ws.Range("C30").Value = 50
$ws.Range("D30").Value = 12907
$ws.Range("E30").Value = 1961
$ws.Range("F30").Value = 128
$ws.Range("G30").Value = 4
$ws.Range("H30").Value = 584

# --- Polonia (row 34) ---
$ws.Range("B34").Value = 12781
$ws.Range("C34").Value = 141
$ws.Range("E34").Value = 9128
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 628

# --- Dinamarca (row 41) ---
$ws.Range("B41").Value = 9158
$ws.Range("C41").Value = 150
$ws.Range("E41").Value = 2349

# --- Filipinas (row 43) ---
$ws.Range("B43").Value = 8488
$ws.Range("C43").Value = 276
$ws.Range("D43").Value = 1043
$ws.Range("E43").Value = 6877
$ws.Range("G43").Value = 10
$ws.Range("H43").Value = 568

# --- Laos (row 185) ---
$ws.Range("D185").Value = 8
$ws.Range("E185").Value = 11

# --- Reorder Santo Tome y Principe / Surinam / Papua Nueva Guinea ---
# Santo Tome y Principe moves up into row 205 with updated counts; Surinam
# and Papua Nueva Guinea each shift down one row, keeping their own figures.
$ws.Range("A205").Value = "Santo Tome y Principe"
$ws.Range("B205").Value = 10
$ws.Range("C205").Value = 2
$ws.Range("D205").Value = 4
$ws.Range("E205").Value = 6
$ws.Range("F205").Value = 0
$ws.Range("G205").Value = 0
$ws.Range("H205").Value = 0

$ws.Range("A206").Value = "Surinam"
$ws.Range("B206").Value = 10
$ws.Range("C206").Value = 0
$ws.Range("D206").Value = 8
$ws.Range("E206").Value = 1
$ws.Range("F206").Value = 0
$ws.Range("G206").Value = 0
$ws.Range("H206").Value = 1

$ws.Range("A207").Value = "Papua Nueva Guinea"
$ws.Range("B207").Value = 8
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 0
$ws.Range("E207").Value = 8
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

$wb.Save()
